# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to several worksheets
# (columns H-N) as captured by the source diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 113.23077
$ws.Range("I5").Value = 110.2
$ws.Range("J5").Value = 123.333336
$ws.Range("K5").Value = 110.2
$ws.Range("L5").Value = 123.333336
$ws.Range("M5").Value = 4.799999999999997
$ws.Range("N5").Value = -353.333336
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = -440
$ws.Range("H39").Value = 243.33333
$ws.Range("I39").Value = 232
$ws.Range("K39").Value = 696
$ws.Range("M39").Value = -400
$ws.Range("H40").Value = 37038964
$ws.Range("J40").Value = 43480304
$ws.Range("L40").Value = 43480304
$ws.Range("N40").Value = -43480654
$ws.Range("H48").Value = 7021.125
$ws.Range("J48").Value = 7021.125
$ws.Range("L48").Value = 21063.375
$ws.Range("N48").Value = -21647.375
$ws.Range("H56").Value = 7021.125
$ws.Range("J56").Value = 7021.125
$ws.Range("L56").Value = 21063.375
$ws.Range("N56").Value = -22131.375
$ws.Range("H116").Value = 4833.6665
$ws.Range("I116").Value = 5050.5
$ws.Range("J116").Value = 4400
$ws.Range("K116").Value = 5050.5
$ws.Range("L116").Value = 4400
$ws.Range("M116").Value = -1608.5
$ws.Range("N116").Value = -11284
$ws.Range("H138").Value = 1231.9
$ws.Range("I138").Value = 649.0294
$ws.Range("J138").Value = 2470.5
$ws.Range("K138").Value = 1947.0882
$ws.Range("L138").Value = 7411.5
$ws.Range("M138").Value = 3192.9118
$ws.Range("N138").Value = -17691.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25550.135
$ws.Range("I32").Value = 26579.352
$ws.Range("J32").Value = 22771.25
$ws.Range("K32").Value = 26579.352
$ws.Range("L32").Value = 22771.25
$ws.Range("M32").Value = -26292.352
$ws.Range("N32").Value = -23345.25
$ws.Range("H61").Value = 2265
$ws.Range("I61").Value = 1040
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1040
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -828
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 870.127
$ws.Range("I74").Value = 797.31915
$ws.Range("J74").Value = 1084
$ws.Range("K74").Value = 797.31915
$ws.Range("L74").Value = 1084
$ws.Range("M74").Value = 76.68084999999996
$ws.Range("N74").Value = -2832
$ws.Range("H77").Value = 870.127
$ws.Range("I77").Value = 797.31915
$ws.Range("J77").Value = 1084
$ws.Range("K77").Value = 3986.59575
$ws.Range("L77").Value = 5420
$ws.Range("M77").Value = 381.4042499999996
$ws.Range("N77").Value = -14156
$ws.Range("H97").Value = 806.41174
$ws.Range("I97").Value = 750.6923
$ws.Range("J97").Value = 987.5
$ws.Range("K97").Value = 750.6923
$ws.Range("L97").Value = 987.5
$ws.Range("M97").Value = -254.6923
$ws.Range("N97").Value = -1979.5
$ws.Range("H132").Value = 1482.125
$ws.Range("I132").Value = 1026.6154
$ws.Range("J132").Value = 2328.0715
$ws.Range("K132").Value = 3079.8462
$ws.Range("L132").Value = 6984.2145
$ws.Range("M132").Value = -549.8462
$ws.Range("N132").Value = -12044.2145
$ws.Range("H136").Value = 2265
$ws.Range("I136").Value = 1040
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3120
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -570
$ws.Range("N136").Value = -14100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5066.25
$ws.Range("I105").Value = 4765.2666
$ws.Range("J105").Value = 5567.8887
$ws.Range("K105").Value = 4765.2666
$ws.Range("L105").Value = 5567.8887
$ws.Range("M105").Value = -3018.2666
$ws.Range("N105").Value = -9061.8887
$ws.Range("H134").Value = 16498.643
$ws.Range("I134").Value = 1268.4906
$ws.Range("J134").Value = 74155.64
$ws.Range("K134").Value = 3805.4718
$ws.Range("L134").Value = 222466.92
$ws.Range("M134").Value = -1270.4718
$ws.Range("N134").Value = -227536.92

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 410.8
$ws.Range("I22").Value = 401.5
$ws.Range("J22").Value = 424.75
$ws.Range("K22").Value = 401.5
$ws.Range("L22").Value = 424.75
$ws.Range("M22").Value = -51.5
$ws.Range("N22").Value = -1124.75
$ws.Range("H31").Value = 2692.1592
$ws.Range("I31").Value = 2772.7742
$ws.Range("J31").Value = 2499.923
$ws.Range("K31").Value = 2772.7742
$ws.Range("L31").Value = 2499.923
$ws.Range("M31").Value = -2477.7742
$ws.Range("N31").Value = -3089.923
$ws.Range("H34").Value = 2692.1592
$ws.Range("I34").Value = 2772.7742
$ws.Range("J34").Value = 2499.923
$ws.Range("K34").Value = 2772.7742
$ws.Range("L34").Value = 2499.923
$ws.Range("M34").Value = -2570.7742
$ws.Range("N34").Value = -2903.923
$ws.Range("H58").Value = 5443.893
$ws.Range("I58").Value = 1531.6316
$ws.Range("J58").Value = 13703.111
$ws.Range("K58").Value = 1531.6316
$ws.Range("L58").Value = 13703.111
$ws.Range("M58").Value = -1328.6316
$ws.Range("N58").Value = -14109.111
$ws.Range("H109").Value = 22074
$ws.Range("J109").Value = 22074
$ws.Range("L109").Value = 22074
$ws.Range("N109").Value = -24154
$ws.Range("H134").Value = 1187.1041
$ws.Range("I134").Value = 965.1667
$ws.Range("K134").Value = 2895.5001
$ws.Range("M134").Value = -360.5001000000002
$ws.Range("H136").Value = 5443.893
$ws.Range("I136").Value = 1531.6316
$ws.Range("J136").Value = 13703.111
$ws.Range("K136").Value = 4594.8948
$ws.Range("L136").Value = 41109.333
$ws.Range("M136").Value = -2044.8948
$ws.Range("N136").Value = -46209.333
$ws.Range("H140").Value = 65061.43
$ws.Range("J140").Value = 65061.43
$ws.Range("L140").Value = 65061.43
$ws.Range("N140").Value = -75421.42999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 190
$ws.Range("I33").Value = 195
$ws.Range("J33").Value = 180
$ws.Range("K33").Value = 1170
$ws.Range("L33").Value = 1080
$ws.Range("M33").Value = -887
$ws.Range("N33").Value = -1646
$ws.Range("H34").Value = 1295.091
$ws.Range("J34").Value = 2580
$ws.Range("L34").Value = 7740
$ws.Range("N34").Value = -7908
$ws.Range("H39").Value = 5586.3125
$ws.Range("J39").Value = 5586.3125
$ws.Range("L39").Value = 16758.9375
$ws.Range("N39").Value = -17346.9375
$ws.Range("H44").Value = 12689
$ws.Range("I44").Value = 252
$ws.Range("J44").Value = 50000
$ws.Range("K44").Value = 756
$ws.Range("L44").Value = 150000
$ws.Range("M44").Value = -358
$ws.Range("N44").Value = -150796
$ws.Range("H46").Value = 2488.7273
$ws.Range("I46").Value = 1980.4
$ws.Range("J46").Value = 2638.2354
$ws.Range("K46").Value = 5941.200000000001
$ws.Range("L46").Value = 7914.706200000001
$ws.Range("M46").Value = -5850.200000000001
$ws.Range("N46").Value = -8096.706200000001
$ws.Range("H55").Value = 27478.842
$ws.Range("I55").Value = 250474.75
$ws.Range("J55").Value = 1244.0294
$ws.Range("K55").Value = 751424.25
$ws.Range("L55").Value = 3732.0882
$ws.Range("M55").Value = -751247.25
$ws.Range("N55").Value = -4086.0882
$ws.Range("H58").Value = 4242.857
$ws.Range("J58").Value = 4242.857
$ws.Range("L58").Value = 12728.571
$ws.Range("N58").Value = -12984.571
$ws.Range("H64").Value = 2285.5715
$ws.Range("I64").Value = 999.5
$ws.Range("K64").Value = 2998.5
$ws.Range("M64").Value = -2728.5
$ws.Range("H67").Value = 2285.5715
$ws.Range("I67").Value = 999.5
$ws.Range("K67").Value = 2998.5
$ws.Range("M67").Value = -2062.5
$ws.Range("H70").Value = 4703.2856
$ws.Range("I70").Value = 3330.75
$ws.Range("K70").Value = 9992.25
$ws.Range("M70").Value = -9677.25
$ws.Range("H73").Value = 4703.2856
$ws.Range("I73").Value = 3330.75
$ws.Range("K73").Value = 9992.25
$ws.Range("M73").Value = -8900.25
$ws.Range("H94").Value = 101984.8
$ws.Range("I94").Value = 125606
$ws.Range("J94").Value = 7500
$ws.Range("K94").Value = 376818
$ws.Range("L94").Value = 22500
$ws.Range("M94").Value = -376142
$ws.Range("N94").Value = -23852
$ws.Range("H100").Value = 4985.4287
$ws.Range("J100").Value = 4985.4287
$ws.Range("L100").Value = 14956.2861
$ws.Range("N100").Value = -16578.2861
$ws.Range("H107").Value = 598928.4399999999
$ws.Range("I107").Value = 993
$ws.Range("J107").Value = 1111444.6
$ws.Range("K107").Value = 2979
$ws.Range("L107").Value = 3334333.8
$ws.Range("M107").Value = -1059
$ws.Range("N107").Value = -3338173.8
$ws.Range("H131").Value = 8231215.5
$ws.Range("I131").Value = 71573270
$ws.Range("J131").Value = 20207.871
$ws.Range("K131").Value = 214719810
$ws.Range("L131").Value = 60623.613
$ws.Range("M131").Value = -214714770
$ws.Range("N131").Value = -70703.613

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 82.72727
$ws.Range("I2").Value = 84.71429000000001
$ws.Range("K2").Value = 84.71429000000001
$ws.Range("M2").Value = 28.28570999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1284.5883
$ws.Range("I46").Value = 1369
$ws.Range("J46").Value = 1164
$ws.Range("K46").Value = 1369
$ws.Range("L46").Value = 1164
$ws.Range("M46").Value = -1181
$ws.Range("N46").Value = -1540
$ws.Range("H93").Value = 1662.3334
$ws.Range("I93").Value = 896.63635
$ws.Range("J93").Value = 2310.2307
$ws.Range("K93").Value = 896.63635
$ws.Range("L93").Value = 2310.2307
$ws.Range("M93").Value = 351.36365
$ws.Range("N93").Value = -4806.2307
$ws.Range("H132").Value = 2098.5667
$ws.Range("I132").Value = 2010.7084
$ws.Range("J132").Value = 2450
$ws.Range("K132").Value = 6032.1252
$ws.Range("L132").Value = 7350
$ws.Range("M132").Value = -3502.1252
$ws.Range("N132").Value = -12410

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 14355.556
$ws.Range("J109").Value = 14355.556
$ws.Range("L109").Value = 14355.556
$ws.Range("N109").Value = -17129.556
$ws.Range("H122").Value = 1079.3334
$ws.Range("I122").Value = 1051
$ws.Range("K122").Value = 3153
$ws.Range("M122").Value = -703
$ws.Range("H132").Value = 907.3090999999999
$ws.Range("I132").Value = 613.1556
$ws.Range("J132").Value = 2231
$ws.Range("K132").Value = 1839.4668
$ws.Range("L132").Value = 6693
$ws.Range("M132").Value = 690.5331999999999
$ws.Range("N132").Value = -11753
$ws.Range("H136").Value = 421.72223
$ws.Range("I136").Value = 257.33334
$ws.Range("K136").Value = 772.0000200000001
$ws.Range("M136").Value = 1777.99998

